$wb = $excel.ActiveWorkbook

# xlNone = -4142
$xlNone = -4142
# xlContinuous = 1
$xlContinuous = 1
# xlPasteFormats = -4122
$xlPasteFormats = -4122
# Border edge indices: xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
$xlEdgeLeft = 7
$xlEdgeRight = 10

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1 becomes the "middle of merged header" cell (top+bottom border only) and
# D1 becomes the "end of merged header" cell (top+bottom+right border), both
# reverting from the bold/centered header font to the plain default font -
# this matches the two new cellXfs entries added by the edit.
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = $xlContinuous
$c1.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
$c1.Borders.Item($xlEdgeRight).LineStyle = $xlNone

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.LineStyle = $xlContinuous
$d1.Borders.Item($xlEdgeLeft).LineStyle = $xlNone

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

# F1/G1 need the exact same two formats as C1/D1 on this sheet - copy the
# already-built formats across (instead of rebuilding borders from scratch)
# so no extra/duplicate style entries get created.
$c1.Copy() | Out-Null
$ws2.Range("C1").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Range("F1").PasteSpecial($xlPasteFormats) | Out-Null

$d1.Copy() | Out-Null
$ws2.Range("D1").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell entirely.
$ws2.Range("G5").ClearContents()
